$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "260.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.93%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.05"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.65%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.702"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.38%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06178"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.69%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.76%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8502"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.70%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9138"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.06%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.34%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04651"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "8.90%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07083"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.91%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03112"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.53%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09044"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.75%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001539"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.59%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006167"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.94%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006048"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.63%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.453"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.03%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.33%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.88%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.085"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.05%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04240"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.38%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.08%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.65%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.11%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-7.80%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03880"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.49%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1112"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.13%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004090"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "7.94%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "9.28%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-10.05%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.15%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.12%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1682"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-23.92%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.12%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.12%"
